$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles/borders/fonts) of row 10 down into the new row 11
$ws.Range("A10:L10").Copy() | Out-Null
$ws.Range("A11:L11").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Set the row height for the newly added row
$ws.Rows.Item(11).RowHeight = 30

# Populate the new test case row
$ws.Range("A11").Value = "CP_AUTO_010"
$ws.Range("B11").Value = "Dispositivos conectados"
$ws.Range("C11").Value = "Positivo"
$ws.Range("D11").Value = "eCenter"
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = ""
$ws.Range("G11").Value = ""
$ws.Range("H11").Value = ""
$ws.Range("I11").Value = ""
$ws.Range("J11").Value = "OK"
$ws.Range("K11").Value = "SI"
$ws.Range("L11").Value = "N/A"

# Update the active selection as reflected in the saved view state
$ws.Range("F8").Select() | Out-Null
